# Update cryptos list - refresh prices / 1h volume percentages, and
# fix the ordering of FraxShare / MXToken rows (43 and 44 swapped places).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.989.16"
$ws.Range("E2").Value = "  +1.01%  "

$ws.Range("D3").Value = "1.555.12"
$ws.Range("E3").Value = "  +0.63%  "

$ws.Range("E4").Value = "  +0.45%  "

$ws.Range("D5").Value = "207.27"
$ws.Range("E5").Value = "  +0.72%  "

$ws.Range("E6").Value = "  +0.97%  "

$ws.Range("E8").Value = "  +1.47%  "

$ws.Range("D9").Value = "21.62"
$ws.Range("E9").Value = "  +1.23%  "

$ws.Range("E10").Value = "  +1.41%  "

$ws.Range("E11").Value = "  +1.07%  "

$ws.Range("D12").Value = "1.776.90"
$ws.Range("E12").Value = "  +0.63%  "

$ws.Range("D13").Value = "1.554.95"
$ws.Range("E13").Value = "  +0.55%  "

$ws.Range("E14").Value = "  +1.42%  "

$ws.Range("D15").Value = "0.515"
$ws.Range("E15").Value = "  +1.00%  "

$ws.Range("D16").Value = "61.92"
$ws.Range("E16").Value = "  +1.28%  "

$ws.Range("D17").Value = "26.971.12"
$ws.Range("E17").Value = "  +0.99%  "

$ws.Range("D18").Value = "215.63"
$ws.Range("E18").Value = "  +1.55%  "

$ws.Range("E19").Value = "  +0.09%  "

$ws.Range("D20").Value = "7.27"
$ws.Range("E20").Value = "  +0.80%  "

$ws.Range("E21").Value = "  +0.44%  "

$ws.Range("E22").Value = "  -0.77%  "

$ws.Range("E23").Value = "  +2.95%  "

$ws.Range("E24").Value = "  -0.47%  "

$ws.Range("D25").Value = "152.37"
$ws.Range("E25").Value = "  -0.18%  "

$ws.Range("E26").Value = "  +2.48%  "

$ws.Range("D27").Value = "14.90"
$ws.Range("E27").Value = "  +0.28%  "

$ws.Range("E28").Value = "  +0.46%  "

$ws.Range("E29").Value = "  +1.42%  "

$ws.Range("E31").Value = "  -0.16%  "

$ws.Range("E32").Value = "  +1.01%  "

$ws.Range("D33").Value = "1.402.41"
$ws.Range("E33").Value = "  +5.04%  "

$ws.Range("E34").Value = "  +3.04%  "

$ws.Range("E35").Value = "  +3.36%  "

$ws.Range("D36").Value = "0.952"
$ws.Range("E36").Value = "  +2.48%  "

$ws.Range("E37").Value = "  +0.44%  "

$ws.Range("E38").Value = "  +0.78%  "

$ws.Range("D39").Value = "0.524"
$ws.Range("E39").Value = "  +0.30%  "

$ws.Range("E40").Value = "  +1.38%  "

$ws.Range("E41").Value = "  +0.47%  "

$ws.Range("D42").Value = "0.989"
$ws.Range("E42").Value = "  -0.64%  "

# Rows 43 and 44 swap: FraxShare moves up to row 43, MXToken moves down to row 44.
$ws.Range("B43").Value = "FraxShare"
$ws.Range("C43").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D43").Value = "5.51"
$ws.Range("E43").Value = "  -3.95%  "

$ws.Range("B44").Value = "MXToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D44").Value = "2.27"
$ws.Range("E44").Value = "  +3.35%  "

$ws.Range("D45").Value = "63.83"
$ws.Range("E45").Value = "  +2.06%  "

$ws.Range("E46").Value = "  +0.24%  "

$ws.Range("D47").Value = "1.690.68"

$ws.Range("D48").Value = "86.28"
$ws.Range("E48").Value = "  +0.47%  "

$ws.Range("E49").Value = "  +1.32%  "

$ws.Range("E50").Value = "  +0.55%  "

$ws.Range("E51").Value = "  +0.49%  "
